$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$data = @(
    @(91.32270166935785, 225451, 334.9940564635958),
    @(91.18905604074934, 38102, 242.687898089172),
    @(87.48757423640602, 147045, 141.3894230769231),
    @(97.12764057113438, 56601, 165.0174927113703),
    @(18.13544486308301, 1980, 14.04255319148936),
    @(31.92446020123002, 228, 20.72727272727273)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 12).Value = $data[$i][0]
    $ws.Cells.Item($row, 13).Value = $data[$i][1]
    $ws.Cells.Item($row, 14).Value = $data[$i][2]
}
